$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column widths: A-C stay the same (30.28515625), D widens, E and F are new
# (ColumnWidth values chosen so the saved stored-width lands as close as
# possible to the target 55.5703125 / 32.85546875 / 15.85546875)
$ws.Columns.Item(4).ColumnWidth = 54.666666666666664
$ws.Columns.Item(5).ColumnWidth = 32.0
$ws.Columns.Item(6).ColumnWidth = 15.0

# New Widget/Signal/Slot rows (10-13)
$ws.Cells.Item(10, 3).Value = "ROI array radio button"
$ws.Cells.Item(10, 4).Value = "Open a dialog window"

$ws.Cells.Item(11, 2).Value = "erase last ROI"

$ws.Cells.Item(12, 2).Value = "create resulting XLSX"

$ws.Cells.Item(13, 2).Value = "processLabel"
$ws.Cells.Item(13, 4).Value = "show the function/method/process which is being executed"

# New columns E:G (class name / method name / line) header
$ws.Cells.Item(1, 5).Value = "class name"
$ws.Cells.Item(1, 6).Value = "method name"
$ws.Cells.Item(1, 7).Value = "line"

# Row 2
$ws.Cells.Item(2, 5).Value = "GUI"
$ws.Cells.Item(2, 6).Value = "motion"
$ws.Cells.Item(2, 7).Value = 1824

# Row 7
$ws.Cells.Item(7, 5).Value = "GUI"
$ws.Cells.Item(7, 6).Value = "slide_images"
$ws.Cells.Item(7, 7).Value = 1848

# Row 8
$ws.Cells.Item(8, 6).Value = "to be added"

# Row 4
$ws.Cells.Item(4, 5).Value = "GUI"
$ws.Cells.Item(4, 6).Value = "draw"
$ws.Cells.Item(4, 7).Value = 2028

# Row 11 (rest of it)
$ws.Cells.Item(11, 5).Value = "GUI"
$ws.Cells.Item(11, 6).Value = "erase_last"
$ws.Cells.Item(11, 7).Value = 2175

# Row 3 (rest of it, reuses existing strings)
$ws.Cells.Item(3, 5).Value = "GUI"
$ws.Cells.Item(3, 6).Value = "motion"

# Match the final selection shown in the diff
$ws.Range("E12").Select()
